$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 4 (row 5) gains a CapacityDone value in column F
$ws.Range("F5").Value = 3

# New sprint rows 6-9 (Sprint 5..Sprint 8)
$ws.Range("A6").Value = "Sprint 5"
$ws.Range("D6").Value = 6

$ws.Range("A7").Value = "Sprint 6"
$ws.Range("D7").Value = 4

$ws.Range("A8").Value = "Sprint 7"
$ws.Range("D8").Value = 3

$ws.Range("A9").Value = "Sprint 8"
$ws.Range("D9").Value = 3

# Copy the date formatting from the existing rows so the new date cells
# reuse the same style (and don't mint a new number format)
$ws.Range("B5:C5").Copy()
$ws.Range("B6:C9").PasteSpecial(-4122)

$ws.Range("B6").Value = 41760
$ws.Range("C6").Value = 41760
$ws.Range("B7").Value = 41761
$ws.Range("C7").Value = 41761
$ws.Range("B8").Value = 41761
$ws.Range("C8").Value = 41761
$ws.Range("B9").Value = 41761
$ws.Range("C9").Value = 41761

$ws.Range("D10").Select()
